$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RemoveCustomerTest")

$ws.Range("A4").Value = "Jose"
$ws.Range("B4").Value = "Silva"
$ws.Range("C4").Value = 654987
$ws.Range("D4").Value = "Customer added successfully"
$ws.Range("E4").Value = "y"

$ws.Activate()
$ws.Range("E5").Select() | Out-Null
